# Applies the "Second order groups with new excel file done!!" edit:
#  - FamilyGroups sheet (UNIFAC-DORTMUND-SurfaceVolume's neighbour tab)
#    gets a new first column "INDEX" (1..23) inserted before the existing
#    "NAME"/"GROUPS" columns, shifting everything right by one column.
#  - The UNIFAC-DORTMUND-SurfaceVolume sheet view scrolls down and loses
#    the "active sheet" flag, because FamilyGroups becomes the active tab.
#  - The FamilyGroups sheet view becomes the active tab with a new
#    selection.

$wb = $excel.ActiveWorkbook

$wsSurface = $wb.Worksheets.Item("UNIFAC-DORTMUND-SurfaceVolume")
$wsFamily  = $wb.Worksheets.Item("FamilyGroups")

# ---------------------------------------------------------------------
# 1. FamilyGroups: insert a new column A ("INDEX") in front of the
#    current A ("NAME") / B ("GROUPS") columns. Columns.Insert shifts
#    every existing cell (values + styles) one column to the right and
#    keeps the dimension/col width metadata in sync.
# ---------------------------------------------------------------------
$wsFamily.Columns.Item(1).Insert()

# Header row (row 1): give the new A1 the same look as the rest of the
# bold/bordered header row (style copied from B1, which used to be A1).
$wsFamily.Range("B1").Copy()
$wsFamily.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$wsFamily.Range("A1").Value = "INDEX"

# Data rows 2..24: number them 1..23, matching the style used by the
# rest of the row (copied from B<n>, which used to be A<n>).
for ($row = 2; $row -le 24; $row++) {
    $srcCell = $wsFamily.Cells.Item($row, 2)
    $dstCell = $wsFamily.Cells.Item($row, 1)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
    $dstCell.Value = $row - 1
}

$wsFamily.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. View state: FamilyGroups becomes the active sheet/tab, with a new
#    selection; UNIFAC-DORTMUND-SurfaceVolume scrolls down (and is no
#    longer the active tab).
# ---------------------------------------------------------------------
$wsSurface.Activate()
$wsSurface.Range("D67").Select()

$wsFamily.Activate()
$wsFamily.Range("O13").Select()
